$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.780.38"
$ws.Range("E2").Value = "  +2.93%  "

# Row 3
$ws.Range("D3").Value = "3.134.48"
$ws.Range("E3").Value = "  +1.88%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.97"
$ws.Range("E5").Value = "  +1.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.36"
$ws.Range("E6").Value = "  +3.54%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "3.128.59"

# Row 9
$ws.Range("E9").Value = "  +1.31%  "

# Row 10
$ws.Range("E10").Value = "  +14.84%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").Value = "  -0.90%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  +0.90%  "

# Row 13
$ws.Range("E13").Value = "  +5.53%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.38"
$ws.Range("E14").Value = "  +5.95%  "

# Row 15
$ws.Range("E15").Value = "  -0.48%  "

# Row 16
$ws.Range("D16").Value = "3.653.67"
$ws.Range("E16").Value = "  +1.89%  "

# Row 17
$ws.Range("E17").Value = "  -1.16%  "

# Row 18
$ws.Range("D18").Value = "63.674.45"
$ws.Range("E18").Value = "  +2.90%  "

# Row 19
$ws.Range("D19").Value = "3.132.78"
$ws.Range("E19").Value = "  +1.79%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.98"
$ws.Range("E20").Value = "  +4.01%  "

# Row 21
$ws.Range("E21").Value = "  +3.30%  "

# Row 22
$ws.Range("E22").Value = "  +0.35%  "

# Row 23
$ws.Range("E23").Value = "  +1.53%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.29"
$ws.Range("E24").Value = "  -3.73%  "

# Row 26
$ws.Range("E26").Value = "  -0.37%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.97"
$ws.Range("E27").Value = "  +8.55%  "

# Row 28
$ws.Range("E28").Value = "  +1.77%  "

# Row 29
$ws.Range("E29").Value = "  -1.72%  "

# Row 30
$ws.Range("E30").Value = "  -0.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.89"
$ws.Range("E31").Value = "  +1.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.22"
$ws.Range("E32").Value = "  +1.28%  "

# Row 33
$ws.Range("E33").Value = "  -2.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0882"
$ws.Range("E34").Value = "  +11.39%  "

# Row 35
$ws.Range("E35").Value = "  +8.36%  "

# Row 36
$ws.Range("E36").Value = "  +1.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.43"
$ws.Range("E37").Value = "  +14.77%  "

# Row 39
$ws.Range("E39").Value = "  +1.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "455.45"
$ws.Range("E40").Value = "  +8.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.76"
$ws.Range("E41").Value = "  -0.32%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0374"
$ws.Range("E42").Value = "  +0.70%  "

# Row 43
$ws.Range("D43").Value = "2.905.70"
$ws.Range("E43").Value = "  -0.76%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.280"
$ws.Range("E44").Value = "  +2.10%  "

# Row 45
$ws.Range("E45").Value = "  +1.81%  "

# Row 46
$ws.Range("E46").Value = "  +2.48%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.90"
$ws.Range("E47").Value = "  +2.85%  "

# Row 48
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.93"
$ws.Range("E48").Value = "  +2.13%  "

# Row 50
$ws.Range("E50").Value = "  +0.49%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.85"
$ws.Range("E51").Value = "  +1.67%  "
